$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the "MCT-3B-Motores CA" class from C3 to C7
$ws.Range("C3").Value = "-"
$ws.Range("C7").Value = "MCT-3B-Motores CA"

# Move the "MEC-3B-Motores CA" class from B4 to B7
$ws.Range("B4").Value = "-"
$ws.Range("B7").Value = "MEC-3B-Motores CA"
